$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Object_Mapping")

# Insert a new row at row 8, shifting existing rows 8-13 down to 9-14
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row with the Wind_onshore entry
$ws.Range("A8").Value = "Wind_onshore"
$ws.Range("B8").Value = "Wind_onshore"
